$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'66.848.47"
$ws.Range("E2").Value = "  -0.40%  "

# Row 3
$ws.Range("D3").Value = "'3.449.17"
$ws.Range("E3").Value = "  -1.66%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'590.26"
$ws.Range("E5").Value = "  -1.05%  "

# Row 6
$ws.Range("D6").Value = "'177.89"
$ws.Range("E6").Value = "  +2.77%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.607"
$ws.Range("E8").Value = "  +4.16%  "

# Row 9
$ws.Range("D9").Value = "'3.448.89"
$ws.Range("E9").Value = "  -1.57%  "

# Row 10
$ws.Range("D10").Value = "'0.136"
$ws.Range("E10").Value = "  +3.09%  "

# Row 11
$ws.Range("D11").Value = "'6.92"
$ws.Range("E11").Value = "  -3.38%  "

# Row 12
$ws.Range("D12").Value = "'0.430"
$ws.Range("E12").Value = "  -0.08%  "

# Row 13
$ws.Range("D13").Value = "'4.052.65"
$ws.Range("E13").Value = "  -1.32%  "

# Row 14
$ws.Range("D14").Value = "'31.47"
$ws.Range("E14").Value = "  +5.28%  "

# Row 15
$ws.Range("E15").Value = "  -0.44%  "

# Row 16
$ws.Range("D16").Value = "'66.901.00"
$ws.Range("E16").Value = "  -0.26%  "

# Row 17
$ws.Range("D17").Value = "'0.0000175"
$ws.Range("E17").Value = "  -2.14%  "

# Row 18
$ws.Range("D18").Value = "'3.456.24"
$ws.Range("E18").Value = "  -1.61%  "

# Row 19
$ws.Range("D19").Value = "'6.22"
$ws.Range("E19").Value = "  -0.95%  "

# Row 20
$ws.Range("D20").Value = "'14.06"
$ws.Range("E20").Value = "  -2.63%  "

# Row 21
$ws.Range("D21").Value = "'386.64"
$ws.Range("E21").Value = "  -1.43%  "

# Row 22
$ws.Range("D22").Value = "'7.88"
$ws.Range("E22").Value = "  -0.84%  "

# Row 23
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.35%  "

# Row 24
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'5.75"
$ws.Range("E24").Value = "  +1.26%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'72.16"
$ws.Range("E25").Value = "  -1.67%  "

# Row 26
$ws.Range("D26").Value = "'0.533"
$ws.Range("E26").Value = "  -0.57%  "

# Row 27
$ws.Range("D27").Value = "'0.0000120"
$ws.Range("E27").Value = "  -1.22%  "

# Row 28
$ws.Range("D28").Value = "'10.24"
$ws.Range("E28").Value = "  +0.59%  "

# Row 29
$ws.Range("D29").Value = "'0.173"
$ws.Range("E29").Value = "  -3.95%  "

# Row 30
$ws.Range("E30").Value = "  +0.47%  "

# Row 31
$ws.Range("D31").Value = "'6.12"
$ws.Range("E31").Value = "  -0.14%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.05"
$ws.Range("E32").Value = "  -0.65%  "

# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.39"
$ws.Range("E33").Value = "  -2.62%  "

# Row 34
$ws.Range("D34").Value = "'23.34"
$ws.Range("E34").Value = "  -1.20%  "

# Row 35
$ws.Range("D35").Value = "'7.30"
$ws.Range("E35").Value = "  -0.88%  "

# Row 36
$ws.Range("E36").Value = "  -0.06%  "

# Row 37
$ws.Range("D37").Value = "'1.58"
$ws.Range("E37").Value = "  -1.93%  "

# Row 38
$ws.Range("D38").Value = "'161.51"
$ws.Range("E38").Value = "  -1.25%  "

# Row 39
$ws.Range("D39").Value = "'0.874"
$ws.Range("E39").Value = "  -0.18%  "

# Row 40
$ws.Range("D40").Value = "'2.77"
$ws.Range("E40").Value = "  +7.58%  "

# Row 41
$ws.Range("D41").Value = "'1.85"
$ws.Range("E41").Value = "  -2.98%  "

# Row 42
$ws.Range("D42").Value = "'6.79"
$ws.Range("E42").Value = "  -0.95%  "

# Row 43
$ws.Range("D43").Value = "'4.65"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("D44").Value = "'25.94"
$ws.Range("E44").Value = "  -0.47%  "

# Row 45
$ws.Range("D45").Value = "'2.769.24"
$ws.Range("E45").Value = "  -1.61%  "

# Row 46
$ws.Range("D46").Value = "'0.0716"
$ws.Range("E46").Value = "  -1.91%  "

# Row 47
$ws.Range("D47").Value = "'26.00"
$ws.Range("E47").Value = "  -4.12%  "

# Row 48
$ws.Range("D48").Value = "'41.00"
$ws.Range("E48").Value = "  -3.45%  "

# Row 49
$ws.Range("D49").Value = "'0.0296"
$ws.Range("E49").Value = "  -2.39%  "

# Row 50
$ws.Range("D50").Value = "'324.98"
$ws.Range("E50").Value = "  -5.33%  "

# Row 51
$ws.Range("E51").Value = "  -3.58%  "
